# Updated cryptos list on Mon Jul 31 03:56:06 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.447.69"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "1.870.28"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7062"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07876"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.92%  "

$ws.Range("D12").Value = "1.895.54"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7046"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.490"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.75%  "

$ws.Range("D17").Value = "29.519.36"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008375"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "256.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.34%  "

$ws.Range("D20").Value = "2.139.25"
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.620"
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1555"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.064"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.499"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.335"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.247"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.209"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05325"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.901"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7481"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.72%  "

$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.91%  "

$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").Value = "1.265.75"
$ws.Range("E39").Value = "  -0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9006"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.959"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.81%  "

$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("E46").Value = "  +1.51%  "

$ws.Range("D47").Value = "2.040.20"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.816"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("E49").Value = "  -0.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.507"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4323"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.31%  "

